# Doing Updates for Financials
# Insert a new "most recent quarter" column before column D on the ENBP sheet,
# shifting the existing D:K data to E:L, then populate the new column D with
# the latest quarter's figures (and re-point the C:L header spans/format).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ENBP")

# 1. Insert a new blank column at D; this shifts D:K -> E:L (values + styles).
$ws.Columns("D:D").Insert()

# 2. The newly inserted column D has no number formatting yet (it inherited
#    column C's format). Copy the formats from the (now shifted) column E so
#    that D matches the date/number styles used by the rest of the table.
#    Only touch the row blocks that actually hold table data, so we don't
#    introduce stray blank cells on the section-heading / spacer rows.
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)

$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)

$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)

# 3. Populate column D with the newest quarter's values.
$newValues = @{
    7  = 43373
    8  = 9300
    9  = "NA"
    10 = "NA"
    12 = "NA"
    13 = 0
    14 = 0
    15 = 0
    17 = 1100
    18 = 8200
    20 = -5200
    21 = 3400
    22 = 0
    23 = 3000
    24 = 400
    25 = 0
    26 = 2600
    27 = 2600
    28 = 0
    29 = "NA"
    30 = 0
    31 = 0
    32 = 5200
    33 = 2600
    34 = 0
    35 = 2600
    38 = 43373
    41 = 17100
    42 = 17500
    43 = 0
    44 = 0
    45 = 0
    46 = 0
    47 = 0
    48 = 25700
    49 = 0
    50 = 0
    51 = 0
    52 = 0
    53 = 0
    54 = 1062200
    57 = 0
    58 = 0
    59 = 0
    60 = 0
    61 = 0
    62 = 0
    63 = 0
    64 = 0
    65 = 0
    66 = 963000
    68 = 0
    69 = 0
    70 = 0
    71 = 0
    72 = 102700
    73 = 0
    74 = 0
    75 = 0
    76 = 99200
    77 = 0
    80 = 43373
    81 = 2600
    83 = 400
    84 = 0
    85 = 0
    86 = 0
    87 = 0
    88 = 0
    89 = 4400
    91 = 0
    92 = 0
    93 = 0
    94 = -30800
    96 = -800
    97 = 0
    98 = 0
    99 = 0
    100 = 10500
    101 = 0
    102 = -15900
}

foreach ($r in $newValues.Keys) {
    $ws.Cells.Item($r, 4).Value2 = $newValues[$r]
}

# Rows 11, 16, 19, 39, 40, 55, 56, 67, 82, 90, 95 stay blank in column D
# (they were blank in every other quarter column too), so nothing else to set.

# Row 91 ("Capital Expenditures") did not just shift with the rest of the
# table - the historical figures for this line were corrected at the same
# time, so set E91:J91 explicitly (K91/L91 keep the values the insert/shift
# already gave them).
$ws.Cells.Item(91, 5).Value2 = 0
$ws.Cells.Item(91, 6).Value2 = 0
$ws.Cells.Item(91, 7).Value2 = 0
$ws.Cells.Item(91, 8).Value2 = 0
$ws.Cells.Item(91, 9).Value2 = -100
$ws.Cells.Item(91, 10).Value2 = 0
